$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data from the latest GitHub Actions run.
# Force text number format on target cells so numeric-looking strings
# (e.g. "215.36", "19.63") are preserved exactly as text, matching the
# original inlineStr cell contents instead of being parsed as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.930.84'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.638.63'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.36'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.63'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.673.55'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.864.75'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.26'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.97'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.935.56'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.08'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.86'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.30'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.45'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.901'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.139.45'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.49'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.38'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.774.57'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.60'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0532'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.66'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.86%  '
